$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.874.38"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.449.73"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "2.456.69"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0984"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").Value = "2.885.60"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "57.819.21"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "2.450.08"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "263.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.37%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0924"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0495"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.07%  "
